$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1907692307692308
$ws.Range("C2").Value = 0.5753846153846154
$ws.Range("J2").Value = 0.03076923076923077
$ws.Range("P2").Value = 0.1076923076923077
$ws.Range("S2").Value = 0.09538461538461539
$ws.Range("B3").Value = 0.01052631578947368
$ws.Range("C3").Value = 0.02631578947368421
$ws.Range("J3").Value = 0.05789473684210526
$ws.Range("P3").Value = 0.7368421052631579
$ws.Range("S3").Value = 0.1684210526315789
$ws.Range("J4").Value = 0.04
$ws.Range("P4").Value = 0.64
$ws.Range("S4").Value = 0.32
$ws.Range("J5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.04329004329004329
$ws.Range("D6").Value = 0.008658008658008658
$ws.Range("F6").Value = 0.03463203463203463
$ws.Range("J6").Value = 0.2640692640692641
$ws.Range("O6").Value = 0.008658008658008658
$ws.Range("Q6").Value = 0.1818181818181818
$ws.Range("R6").Value = 0.08658008658008658
$ws.Range("S6").Value = 0.3722943722943723
$ws.Range("B7").Value = 0.09289617486338798
$ws.Range("D7").Value = 0.0273224043715847
$ws.Range("F7").Value = 0.06010928961748634
$ws.Range("J7").Value = 0.1475409836065574
$ws.Range("O7").Value = 0.01092896174863388
$ws.Range("Q7").Value = 0.180327868852459
$ws.Range("R7").Value = 0.09836065573770492
$ws.Range("S7").Value = 0.3825136612021858
$ws.Range("B8").Value = 0.08523908523908524
$ws.Range("D8").Value = 0.02494802494802495
$ws.Range("F8").Value = 0.06444906444906445
$ws.Range("J8").Value = 0.1185031185031185
$ws.Range("O8").Value = 0.006237006237006237
$ws.Range("Q8").Value = 0.185031185031185
$ws.Range("R8").Value = 0.08108108108108109
$ws.Range("S8").Value = 0.4345114345114345
$ws.Range("B9").Value = 0.08502024291497975
$ws.Range("D9").Value = 0.008097165991902834
$ws.Range("E9").Value = 0.004048582995951417
$ws.Range("F9").Value = 0.04453441295546558
$ws.Range("J9").Value = 0.09716599190283401
$ws.Range("O9").Value = 0.004048582995951417
$ws.Range("Q9").Value = 0.2186234817813765
$ws.Range("R9").Value = 0.07692307692307693
$ws.Range("S9").Value = 0.4615384615384616
$ws.Range("B10").Value = 0.1324655436447167
$ws.Range("D10").Value = 0.02373660030627871
$ws.Range("E10").Value = 0.001531393568147014
$ws.Range("F10").Value = 0.06814701378254211
$ws.Range("J10").Value = 0.1202143950995406
$ws.Range("O10").Value = 0.009188361408882083
$ws.Range("Q10").Value = 0.2297090352220521
$ws.Range("R10").Value = 0.06967840735068913
$ws.Range("S10").Value = 0.3453292496171516
$ws.Range("G11").Value = 0.1528239202657807
$ws.Range("J11").Value = 0.1229235880398671
$ws.Range("K11").Value = 0.2292358803986711
$ws.Range("L11").Value = 0.478405315614618
$ws.Range("S11").Value = 0.01661129568106312
$ws.Range("G12").Value = 0.7302631578947368
$ws.Range("J12").Value = 0.1907894736842105
$ws.Range("L12").Value = 0.05263157894736842
$ws.Range("S12").Value = 0.02631578947368421
$ws.Range("G13").Value = 0.6730769230769231
$ws.Range("J13").Value = 0.2884615384615384
$ws.Range("S13").Value = 0.03846153846153846
$ws.Range("G14").Value = 0.25
$ws.Range("J14").Value = 0.5
$ws.Range("S14").Value = 0.25
$ws.Range("F15").Value = 0.01932367149758454
$ws.Range("H15").Value = 0.1835748792270532
$ws.Range("I15").Value = 0.106280193236715
$ws.Range("J15").Value = 0.3719806763285024
$ws.Range("K15").Value = 0.06763285024154589
$ws.Range("M15").Value = 0.01449275362318841
$ws.Range("O15").Value = 0.06280193236714976
$ws.Range("S15").Value = 0.1739130434782609
$ws.Range("F16").Value = 0.01515151515151515
$ws.Range("H16").Value = 0.1616161616161616
$ws.Range("I16").Value = 0.09595959595959595
$ws.Range("J16").Value = 0.3535353535353535
$ws.Range("K16").Value = 0.1060606060606061
$ws.Range("M16").Value = 0.04040404040404041
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.1717171717171717
$ws.Range("F17").Value = 0.01343570057581574
$ws.Range("H17").Value = 0.1938579654510557
$ws.Range("I17").Value = 0.1132437619961612
$ws.Range("J17").Value = 0.3838771593090211
$ws.Range("K17").Value = 0.07293666026871401
$ws.Range("M17").Value = 0.02111324376199616
$ws.Range("N17").Value = 0.003838771593090211
$ws.Range("O17").Value = 0.0671785028790787
$ws.Range("S17").Value = 0.1305182341650672
$ws.Range("F18").Value = 0.01612903225806452
$ws.Range("H18").Value = 0.1720430107526882
$ws.Range("I18").Value = 0.1075268817204301
$ws.Range("J18").Value = 0.3655913978494624
$ws.Range("K18").Value = 0.1021505376344086
$ws.Range("M18").Value = 0.02150537634408602
$ws.Range("O18").Value = 0.08602150537634409
$ws.Range("S18").Value = 0.1290322580645161
$ws.Range("F19").Value = 0.02190332326283988
$ws.Range("H19").Value = 0.2122356495468278
$ws.Range("I19").Value = 0.09516616314199396
$ws.Range("J19").Value = 0.3595166163141994
$ws.Range("K19").Value = 0.1012084592145015
$ws.Range("M19").Value = 0.02190332326283988
$ws.Range("N19").Value = 0.001510574018126888
$ws.Range("O19").Value = 0.06570996978851963
$ws.Range("S19").Value = 0.1208459214501511
